# Updates cryptos list values (price + 1h volume %) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Preserve original text (inline-string) cell type: force the cell to
    # Text format before assigning, then restore General/Normal so no
    # stray number-format/style is left behind on the cell.
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextCell "D2" "63.693.40"
$ws.Range("E2").Value = "  -1.22%  "

Set-TextCell "D3" "2.638.00"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextCell "D5" "579.44"
$ws.Range("E5").Value = "  +0.10%  "

Set-TextCell "D6" "155.20"
$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -4.29%  "

Set-TextCell "D9" "2.635.18"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("E13").Value = "  +0.89%  "

Set-TextCell "D14" "28.35"
$ws.Range("E14").Value = "  -0.93%  "

Set-TextCell "D15" "3.115.41"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("E16").Value = "  -1.94%  "

Set-TextCell "D17" "63.651.70"
$ws.Range("E17").Value = "  -1.00%  "

Set-TextCell "D18" "2.644.19"
$ws.Range("E18").Value = "  +0.52%  "

Set-TextCell "D19" "12.11"
$ws.Range("E19").Value = "  -1.35%  "

$ws.Range("E20").Value = "  +3.21%  "

Set-TextCell "D22" "344.22"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("E23").Value = "  +0.27%  "

Set-TextCell "D24" "68.00"
$ws.Range("E24").Value = "  +0.09%  "

Set-TextCell "D25" "1.88"
$ws.Range("E25").Value = "  +7.34%  "

$ws.Range("E26").Value = "  -4.15%  "

Set-TextCell "D27" "600.77"
$ws.Range("E27").Value = "  +6.14%  "

Set-TextCell "D28" "9.21"
$ws.Range("E28").Value = "  -1.98%  "

$ws.Range("E29").Value = "  +1.60%  "

Set-TextCell "D30" "8.10"
$ws.Range("E30").Value = "  +2.29%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("E34").Value = "  +0.79%  "

Set-TextCell "D35" "6.55"
$ws.Range("E35").Value = "  -1.63%  "

Set-TextCell "D36" "5.43"
$ws.Range("E36").Value = "  +2.44%  "

Set-TextCell "D37" "0.402"
$ws.Range("E37").Value = "  -2.51%  "

Set-TextCell "D38" "0.999"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("E40").Value = "  -2.27%  "

Set-TextCell "D41" "150.86"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("E42").Value = "  -0.03%  "

Set-TextCell "D43" "2.54"
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("E44").Value = "  -0.65%  "

Set-TextCell "D45" "161.09"
$ws.Range("E45").Value = "  +1.73%  "

Set-TextCell "D46" "24.09"
$ws.Range("E46").Value = "  +3.99%  "

Set-TextCell "D47" "3.89"
$ws.Range("E47").Value = "  -2.42%  "

Set-TextCell "D48" "0.0584"
$ws.Range("E48").Value = "  -2.67%  "

$ws.Range("E49").Value = "  -0.59%  "

Set-TextCell "D50" "0.0999"
$ws.Range("E50").Value = "  -2.47%  "

Set-TextCell "D51" "0.0247"
$ws.Range("E51").Value = "  -1.51%  "
